$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("D1").Value = "output"
$ws.Range("E1").Value = "comment"
$ws.Range("F1").Value = "todo"
$ws.Range("G1").Value = ""

# --- Row 2: move old F2 judge text into D2, add new comment (E2), new F2 ---
$ws.Range("D2").Value = "我帮你看了一下哈～这个广告使用了模糊的用词和夸大的手法，来宣传莎普爱思滴眼液对白内障的治疗效果。虽然滴眼液可能对某些人有帮助，但这种夸大和不清晰的表述可能存在误导性。建议你谨慎对待这样的广告宣传，并且在购买前咨询专业的医生或药剂师，获取更准确的产品信息和使用建议。避免购买🚫"
$ws.Range("E2").Value = "检索到减肥茶和糖尿病的知识，输出完全错误"
$ws.Range("F2").Value = "知识库"
$ws.Range("G2").Value = ""

# --- Row 3: new record (曲美减肥茶) ---
$ws.Range("A3").Value = "曲美减肥茶"
$ws.Range("B3").Value = "“用曲美，好身材”`n“3周重现轻盈体态，6周瘦出完美曲线”"
$ws.Range("C3").Value = "数据造假,夸大失实"
$ws.Range("D3").Value = "我帮你看了下哈～🔍🏋️‍♂️🚗`n我的结论是 👉 可信 👍✨`n【虚假诈骗程度：⭐ 1/10星】`n理由：`n这个广告语强调了产品的效果和速度，但没有明显的夸大或虚假承诺。`n曲美家居在之前的对话中表现出对客户服务的重视，并且提供了真实的产品信息。`n用户在之前的对话中提到了对曲美家居产品的满意，这增加了广告的可信度。`n建议：`n如果对产品感兴趣，可以进一步了解产品的详细信息和用户评价。`n可以考虑咨询专业的健康顾问，以确保选择适合个人需求的产品。`n别担心，咱们一起留个心眼 👀✨"
$ws.Range("E3").Value = "”之前的对话“?没能运用相关知识判断3周重现...的真假"
$ws.Range("F3").Value = "知识库，工具拓展能力"

# --- Row 4: new record (雅诗兰黛) ---
$ws.Range("A4").Value = "雅诗兰黛"
$ws.Range("B4").Value = "“年轻指数+77%”`n“柔润+17% 平滑+20% 透亮+15%”"
$ws.Range("C4").Value = "数据造假,用词模糊"
$ws.Range("D4").Value = "😊✨【我帮你看了下哈～】这些描述看起来像是某种护肤产品的效果宣传对吧？👍`n👉 我的结论是：无法判断`n【虚假诈骗程度：⭐ 1/10星】`n理由：这些描述是正面的产品效果宣传，但未提供具体的产品信息或来源，无法直接判断其真实性。建议提供更多详细信息以便进一步分析。`n别担心，咱们一起留个心眼 👀✨`n如果有任何具体的产品信息或来源，欢迎随时告诉我，我会帮您进一步分析哈～🔍"
$ws.Range("E4").Value = "对于可疑的数据不加怀疑"
$ws.Range("F4").Value = "prompt，知识库"
